$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21: Book and a Hard Place
$ws.Range("H21").Value = 37585
$ws.Range("I21").Value = 37585
$ws.Range("K21").Value = 37585
$ws.Range("M21").Value = -37117

# Row 23: There's Something about Bury
$ws.Range("H23").Value = 37585
$ws.Range("I23").Value = 37585
$ws.Range("K23").Value = 37585
$ws.Range("M23").Value = -37351

# Row 31: Hush Little Wailer
$ws.Range("H31").Value = 3473.25
$ws.Range("J31").Value = 4464.3335
$ws.Range("L31").Value = 13393.0005
$ws.Range("N31").Value = -13853.0005

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 269.7857
$ws.Range("I33").Value = 275.27274
$ws.Range("K33").Value = 275.27274
$ws.Range("M33").Value = -46.27274

# Row 51: A Bile Business
$ws.Range("H51").Value = 7431.625
$ws.Range("I51").Value = 13871.8
$ws.Range("J51").Value = 5736.8423
$ws.Range("K51").Value = 13871.8
$ws.Range("L51").Value = 5736.8423
$ws.Range("M51").Value = -13387.8
$ws.Range("N51").Value = -6704.8423

# Row 64: Forged from the Void
$ws.Range("H64").Value = 3739
$ws.Range("I64").Value = 3739
$ws.Range("K64").Value = 3739
$ws.Range("M64").Value = -3491

# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 3739
$ws.Range("I67").Value = 3739
$ws.Range("K67").Value = 3739
$ws.Range("M67").Value = -2881

# Row 99: Rumor Has It
$ws.Range("H99").Value = 3799.5715
$ws.Range("J99").Value = 4399.5
$ws.Range("L99").Value = 13198.5
$ws.Range("N99").Value = -16194.5

# Row 100: Asking for a Friend
$ws.Range("H100").Value = 6423.231
$ws.Range("I100").Value = 2690.4614
$ws.Range("K100").Value = 2690.4614
$ws.Range("M100").Value = -2149.4614

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 4758.2856
$ws.Range("I137").Value = 3436.3333
$ws.Range("K137").Value = 10308.9999
$ws.Range("M137").Value = -7758.999899999999

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3277.75
$ws.Range("J138").Value = 3287.8572
$ws.Range("L138").Value = 9863.571599999999
$ws.Range("N138").Value = -20143.5716

$ws = $wb.Worksheets.Item("ARM")
# Row 30: Not Enough Headroom
$ws.Range("H30").Value = 28666
$ws.Range("I30").Value = 9
$ws.Range("J30").Value = 42994.5
$ws.Range("K30").Value = 9
$ws.Range("L30").Value = 42994.5
$ws.Range("M30").Value = 141
$ws.Range("N30").Value = -43294.5

$ws = $wb.Worksheets.Item("BSM")
# Row 21: Awl or Nothing
$ws.Range("H21").Value = 67499.5
$ws.Range("J21").Value = 67499.5
$ws.Range("L21").Value = 67499.5
$ws.Range("N21").Value = -67971.5

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 2527
$ws.Range("I58").Value = 2218.5386
$ws.Range("J58").Value = 3028.25
$ws.Range("K58").Value = 2218.5386
$ws.Range("L58").Value = 3028.25
$ws.Range("M58").Value = -2015.5386
$ws.Range("N58").Value = -3434.25

# Row 136: Turali Quality
$ws.Range("H136").Value = 2527
$ws.Range("I136").Value = 2218.5386
$ws.Range("J136").Value = 3028.25
$ws.Range("K136").Value = 6655.6158
$ws.Range("L136").Value = 9084.75
$ws.Range("M136").Value = -4105.6158
$ws.Range("N136").Value = -14184.75

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 1409.4
$ws.Range("J5").Value = 1798.25
$ws.Range("L5").Value = 5394.75
$ws.Range("N5").Value = -5618.75

# Row 14: Keep Your Powder Dry
$ws.Range("H14").Value = 10280.454
$ws.Range("I14").Value = 10280.454
$ws.Range("K14").Value = 30841.362
$ws.Range("M14").Value = -30668.362

# Row 55: Pagan Pastries
$ws.Range("H55").Value = 5921.2
$ws.Range("I55").Value = 2098.6
$ws.Range("K55").Value = 6295.799999999999
$ws.Range("M55").Value = -6118.799999999999

# Row 56: Culture Club
$ws.Range("H56").Value = 15312.823
$ws.Range("I56").Value = 15312.823
$ws.Range("K56").Value = 15312.823
$ws.Range("M56").Value = -14782.823

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 4508.778
$ws.Range("J131").Value = 10704.75
$ws.Range("L131").Value = 32114.25
$ws.Range("N131").Value = -42194.25

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 1409.4
$ws.Range("J135").Value = 1798.25
$ws.Range("L135").Value = 16184.25
$ws.Range("N135").Value = -21254.25

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers
$ws.Range("H2").Value = 297.125
$ws.Range("I2").Value = 75.59999999999999
$ws.Range("J2").Value = 666.3333
$ws.Range("K2").Value = 75.59999999999999
$ws.Range("L2").Value = 666.3333
$ws.Range("M2").Value = 37.40000000000001
$ws.Range("N2").Value = -892.3333

# Row 18: Gorgeous Gorget
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 288
$ws.Range("N18").ClearContents()

# Row 55: If You've Got It, Flaunt It
$ws.Range("H55").Value = 14999
$ws.Range("I55").Value = 4999
$ws.Range("K55").Value = 4999
$ws.Range("M55").Value = -4672

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 8802430
$ws.Range("I22").Value = 44010068
$ws.Range("J22").Value = 521.5
$ws.Range("K22").Value = 44010068
$ws.Range("L22").Value = 521.5
$ws.Range("M22").Value = -44009773
$ws.Range("N22").Value = -1111.5

# Row 27: Fire and Hide
$ws.Range("H27").Value = 8802430
$ws.Range("I27").Value = 44010068
$ws.Range("J27").Value = 521.5
$ws.Range("K27").Value = 44010068
$ws.Range("L27").Value = 521.5
$ws.Range("M27").Value = -44009961
$ws.Range("N27").Value = -735.5

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 1230
$ws.Range("I55").Value = 731
$ws.Range("K55").Value = 731
$ws.Range("M55").Value = -558

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 6947625.5
$ws.Range("I68").Value = 20835226
$ws.Range("J68").Value = 3825
$ws.Range("K68").Value = 20835226
$ws.Range("L68").Value = 3825
$ws.Range("M68").Value = -20834477
$ws.Range("N68").Value = -5323

# Row 69: Maybe He's a Lion
$ws.Range("H69").Value = 99999
$ws.Range("J69").Value = 99999
$ws.Range("L69").Value = 99999
$ws.Range("N69").Value = -101621

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 6947625.5
$ws.Range("I71").Value = 20835226
$ws.Range("J71").Value = 3825
$ws.Range("K71").Value = 104176130
$ws.Range("L71").Value = 19125
$ws.Range("M71").Value = -104172386
$ws.Range("N71").Value = -26613

# Row 72: The Wyvern of It (L)
$ws.Range("H72").Value = 99999
$ws.Range("J72").Value = 99999
$ws.Range("L72").Value = 299997
$ws.Range("N72").Value = -308109

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 2787.8572
$ws.Range("I136").Value = 2752.8333
$ws.Range("K136").Value = 8258.499899999999
$ws.Range("M136").Value = -5708.499899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 11877.625
$ws.Range("I96").Value = 8185.875
$ws.Range("J96").Value = 15569.375
$ws.Range("K96").Value = 8185.875
$ws.Range("L96").Value = 15569.375
$ws.Range("M96").Value = -6812.875
$ws.Range("N96").Value = -18315.375
